# "Last modif for Dashboard V04"
# Highlight the "Compare 2 or all microrobot using best and bad case" bullet
# (4th paragraph of the "ZoneTexte 14" textbox on slide 1) in yellow, matching
# the highlight already applied to the surrounding bullet points.

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(1)
$sh = $s.Shapes.Item(12)          # "ZoneTexte 14"
$tr = $sh.TextFrame.TextRange

$para = $tr.Paragraphs(4, 1)      # "Compare 2 or all microrobot using best and bad case"
$para.Font.Highlight = 65535      # yellow (RGB 255,255,0) -> srgbClr val="FFFF00"
